$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Move the "Rounded Rectangle 3" card up slightly (a:off y: 130629 -> 121921 EMU)
$roundedRect = $s.Shapes.Item(2)
$roundedRect.Top = 121921 / 12700

# Fill in the "Username or Email" input field ("Rectangle 7") with the name "John"
$usernameField = $s.Shapes.Item(6)
$usernameField.TextFrame.TextRange.Text = "John"
